$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = -3
$ws.Range("F6").Value = -2
$ws.Range("F7").Value = -4
$ws.Range("F8").Value = -1
$ws.Range("F9").Value = -1
$ws.Range("F11").Value = -1
$ws.Range("F12").Value = -4
$ws.Range("F13").Value = -4
$ws.Range("F17").Value = 7
$ws.Range("F19").Value = 6
$ws.Range("F20").Value = 0
